$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text, $styleSourceAddr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Header text updates
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# Numeric value updates
$ws.Range("N15").Value = -77.777777777777
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -53.846153846153
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = -34.146341463414
$ws.Range("L16").Value = -10
$ws.Range("M16").Value = -37.209302325581
$ws.Range("N16").Value = -88.842975206611
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -26.666666666666
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 62
$ws.Range("K17").Value = -25.806451612903
$ws.Range("L17").Value = 43.75
$ws.Range("M17").Value = 411.111111111111
$ws.Range("N17").Value = -11.538461538461
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -26.666666666666
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 72
$ws.Range("K18").Value = -47.222222222222
$ws.Range("L18").Value = 2.702702702702
$ws.Range("M18").Value = -24
$ws.Range("N18").Value = -81.818181818181
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 44
$ws.Range("H19").Value = 12.820512820512
$ws.Range("I19").Value = 128
$ws.Range("J19").Value = 172
$ws.Range("K19").Value = -25.581395348837
$ws.Range("L19").Value = -15.789473684210
$ws.Range("M19").Value = 4.918032786885
$ws.Range("N19").Value = 0.787401574803
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 800
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 142.857142857143
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 60
$ws.Range("L20").Value = 152.631578947368
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -86.127167630057
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 27.777777777777
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 291
$ws.Range("J21").Value = 380
$ws.Range("K21").Value = -23.421052631578
$ws.Range("L21").Value = 5.434782608695
$ws.Range("M21").Value = 16.4
$ws.Range("N21").Value = -70.516717325228
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("M22").Value = -10
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -3.703703703703
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = -19.166666666666
$ws.Range("I24").Value = 329
$ws.Range("J24").Value = 486
$ws.Range("K24").Value = -32.304526748971
$ws.Range("L24").Value = -4.637681159420
$ws.Range("M24").Value = 38.235294117647
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 5.555555555555
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 89
$ws.Range("H25").Value = -13.483146067415
$ws.Range("I25").Value = 219
$ws.Range("J25").Value = 384
$ws.Range("K25").Value = -42.96875
$ws.Range("L25").Value = -23.693379790940
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 5.555555555555
$ws.Range("I26").Value = 65
$ws.Range("J26").Value = 94
$ws.Range("K26").Value = -30.851063829787
$ws.Range("L26").Value = -16.666666666666
$ws.Range("M26").Value = 3.174603174603
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 7
$ws.Range("I28").Value = 29
$ws.Range("K28").Value = 163.636363636364
$ws.Range("L28").Value = 107.142857142857
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 2

# Text (string) value updates with style matching
Set-TextCell "D14" "0" "C14"
Set-TextCell "E14" "***.*" "M14"
Set-TextCell "D22" "0" "C22"
Set-TextCell "E22" "***.*" "N22"
Set-TextCell "C23" "0" "D15"
Set-TextCell "D23" "0" "D15"
Set-TextCell "E23" "***.*" "N23"
Set-TextCell "C27" "0" "D29"
Set-TextCell "D27" "0" "D29"
Set-TextCell "E27" "***.*" "M27"
Set-TextCell "G28" "0" "D28"
Set-TextCell "H28" "***.*" "E28"
Set-TextCell "C29" "0" "D29"
Set-TextCell "C30" "0" "D30"
